$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 18) appended after the existing last row (row 17),
# containing the corrected positive phrase percentage value.
$row = 18

$ws.Cells.Item($row, 1).Value = (Get-Date -Year 2016 -Month 9 -Day 11 -Hour 14 -Minute 52 -Second 24)
$ws.Cells.Item($row, 2).Value = 34
$ws.Cells.Item($row, 3).Value = 53
$ws.Cells.Item($row, 4).Value = 45
$ws.Cells.Item($row, 5).Value = 53
$ws.Cells.Item($row, 6).Value = 6
$ws.Cells.Item($row, 7).Value = 5964
$ws.Cells.Item($row, 8).Value = 8713
$ws.Cells.Item($row, 9).Value = 1072
$ws.Cells.Item($row, 10).Value = 116
$ws.Cells.Item($row, 11).Value = 99
$ws.Cells.Item($row, 12).Value = 31
$ws.Cells.Item($row, 13).Value = 2
$ws.Cells.Item($row, 14).Value = "Bag"
